# HFDA - ch3part2: Add constraints on sales_volume
#
# This reproduces what happens when the Solver "Add Constraint" dialog is
# used twice more on the Bathing Friends Unlimited model (Sheet1):
#   - solver_lhs3 is repointed from B6 -> B5 (existing 2nd constraint slot
#     reused for the B5 side of a <= / >= pair)
#   - two brand-new constraint slots (index 4 and 5) are appended, each
#     carrying its own lhs/rel/rhs triplet
#   - solver_num goes from 3 -> 5 constraints
#   - the rhs values shuffle: rhs2 150, rhs3 400 (was rhs2), new rhs4 300
#     (was rhs3), new rhs5 50
# plus the underlying input cells (B5, B6) that drive the model change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two model input cells -------------------------------------
$ws.Range('B5').Value = 150
$ws.Range('B6').Value = 50

# --- Update Solver constraint defined names -------------------------------

# solver_lhs3 now refers to B5 instead of B6
$wb.Names.Item('solver_lhs3').RefersTo = '=Sheet1!$B$5'

# New constraint slots 4 and 5 (sheet-scoped + hidden, matching the rest of
# the solver_* names)
$n = $ws.Names.Add('solver_lhs4', '=Sheet1!$B$6')
$n.Visible = $false

$n = $ws.Names.Add('solver_lhs5', '=Sheet1!$B$6')
$n.Visible = $false

$n = $ws.Names.Add('solver_rel4', '=1')
$n.Visible = $false

$n = $ws.Names.Add('solver_rel5', '=1')
$n.Visible = $false

$n = $ws.Names.Add('solver_rhs4', '=300')
$n.Visible = $false

$n = $ws.Names.Add('solver_rhs5', '=50')
$n.Visible = $false

# Existing rhs2 / rhs3 values shift as part of the renumbering
$wb.Names.Item('solver_rhs2').RefersTo = '=150'
$wb.Names.Item('solver_rhs3').RefersTo = '=400'

# Total constraint count
$wb.Names.Item('solver_num').RefersTo = '=5'

# --- Match the saved selection --------------------------------------------
$ws.Range('K12').Select()
